$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($row, $col, $text) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

Set-TextValue 2 4 "42.751.25"
Set-TextValue 2 5 "  -0.48%  "
Set-TextValue 3 4 "2.545.00"
Set-TextValue 3 5 "  +0.13%  "
Set-TextValue 4 5 "  +0.10%  "
Set-TextValue 5 4 "310.03"
Set-TextValue 5 5 "  -2.59%  "
Set-TextValue 6 4 "99.14"
Set-TextValue 6 5 "  +1.66%  "
Set-TextValue 7 5 "  -0.90%  "
Set-TextValue 8 5 "  +0.02%  "
Set-TextValue 9 5 "  -0.73%  "
Set-TextValue 10 4 "35.91"
Set-TextValue 10 5 "  -1.24%  "
Set-TextValue 11 5 "  -1.63%  "
Set-TextValue 12 4 "7.40"
Set-TextValue 12 5 "  -2.75%  "
Set-TextValue 13 5 "  -1.18%  "
Set-TextValue 14 4 "2.937.88"
Set-TextValue 14 5 "  +0.26%  "
Set-TextValue 15 4 "15.90"
Set-TextValue 15 5 "  +4.63%  "
Set-TextValue 16 4 "2.573.61"
Set-TextValue 16 5 "  +4.05%  "
Set-TextValue 17 4 "0.839"
Set-TextValue 17 5 "  -1.76%  "
Set-TextValue 18 4 "42.769.65"
Set-TextValue 18 5 "  -0.57%  "
Set-TextValue 19 5 "  -1.91%  "
Set-TextValue 20 4 "12.38"
Set-TextValue 20 5 "  -3.65%  "
Set-TextValue 21 5 "  -1.68%  "
Set-TextValue 22 4 "69.43"
Set-TextValue 23 4 "247.79"
Set-TextValue 23 5 "  -2.84%  "
Set-TextValue 24 5 "  -1.90%  "
Set-TextValue 25 5 "  -0.08%  "
Set-TextValue 26 4 "26.60"
Set-TextValue 26 5 "  +0.09%  "
Set-TextValue 27 5 "  +0.00%  "
Set-TextValue 28 4 "2.36"
Set-TextValue 28 5 "  -2.13%  "
Set-TextValue 29 4 "40.03"
Set-TextValue 29 5 "  -1.81%  "
Set-TextValue 30 4 "10.11"
Set-TextValue 30 5 "  -3.51%  "
Set-TextValue 31 4 "159.08"
Set-TextValue 31 5 "  +0.47%  "
Set-TextValue 32 4 "5.73"
Set-TextValue 32 5 "  -3.36%  "
Set-TextValue 33 5 "  +0.91%  "
Set-TextValue 34 4 "3.29"
Set-TextValue 34 5 "  -2.35%  "
Set-TextValue 35 4 "2.08"
Set-TextValue 35 5 "  -4.03%  "
Set-TextValue 36 5 "  -3.43%  "
Set-TextValue 37 5 "  +5.26%  "
Set-TextValue 38 5 "  -3.66%  "
Set-TextValue 39 5 "  -1.31%  "
Set-TextValue 40 5 "  -0.85%  "
Set-TextValue 41 4 "22.53"
Set-TextValue 41 5 "  +1.41%  "
Set-TextValue 42 4 "4.14"
Set-TextValue 42 5 "  +7.65%  "
Set-TextValue 43 5 "  -0.11%  "
Set-TextValue 44 5 "  -1.44%  "
Set-TextValue 45 5 "  -1.64%  "
Set-TextValue 46 4 "1.989.52"
Set-TextValue 46 5 "  -1.42%  "
Set-TextValue 47 4 "9.05"
Set-TextValue 47 5 "  -1.23%  "
Set-TextValue 48 4 "2.785.35"
Set-TextValue 48 5 "  -0.06%  "
Set-TextValue 49 4 "81.27"
Set-TextValue 49 5 "  -4.03%  "
Set-TextValue 50 5 "  +0.38%  "
Set-TextValue 51 4 "73.41"
Set-TextValue 51 5 "  -4.14%  "
